# Updated capital structure database
# Refresh the Denmark / Insurance (General) rows (rows 2 and 3) with the
# latest database figures. The expected_growth_eps_next_5_years column (F)
# is no longer populated for this dataset, so its value is cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = [ordered]@{
    "D"  = 0.0793
    "E"  = -0.0332
    "G"  = 0.08899777790439291
    "H"  = 0.08899777790439291
    "I"  = 0.0611930944105749
    "J"  = 0.0473104819472803
    "K"  = 155.4
    "L"  = 0.04427098171044385
    "M"  = 116.9
    "N"  = 0.03085409628378378
    "O"  = 0.7522522522522522
    "P"  = 116.9
    "Q"  = 0.03085409628378378
    "R"  = 0.7522522522522522
    "U"  = 550.4
    "V"  = 0.1452702702702703
    "W"  = 0.1790116346043083
    "X"  = 0.04315583832409947
    "Y"  = 0.1358557962802088
    "Z"  = 3.990677580718508
    "AA" = 0.1888008796399992
    "AB" = 0.04037935754891815
    "AC" = 0.148421522091081
    "AD" = 414.1
    "AF" = 414.1
    "AG" = -136.3
    "AH" = 0.09852720740441123
    "AI" = 0.2943978387601309
    "AJ" = -0.03731690622861053
    "AK" = -0.1591917776220509
    "AL" = 15.9
    "AM" = 15.9
    "AN" = 1.881417537482962
    "AO" = 13.50943396226415
    "AP" = -0.6192639709223079
    "AQ" = 13.50943396226415
}

foreach ($row in 2, 3) {
    # Column F (expected_growth_eps_next_5_years) no longer has a value.
    $ws.Range("F$row").Value = $null

    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
